$d = $word.ActiveDocument

$replacements = @(
    @("52×93=4836", "97×73=7081"),
    @("78×11=858", "72×13=936"),
    @("52×31=1612", "42×72=3024"),
    @("42×82=3444", "88×30=2640"),
    @("46×49=2254", "15×89=1335"),
    @("86×63=5418", "76×15=1140"),
    @("93×49=4557", "22×71=1562"),
    @("39×62=2418", "11×18=198"),
    @("52×72=3744", "83×16=1328"),
    @("57×98=5586", "93×25=2325"),
    @("19×86=1634", "81×60=4860"),
    @("55×84=4620", "66×61=4026"),
    @("66×38=2508", "32×56=1792"),
    @("95×81=7695", "66×88=5808"),
    @("48×89=4272", "60×91=5460"),
    @("28×15=420", "74×64=4736"),
    @("46×26=1196", "66×71=4686"),
    @("74×12=888", "77×26=2002"),
    @("65×12=780", "49×52=2548"),
    @("31×22=682", "60×45=2700"),
    @("52×18=936", "46×46=2116"),
    @("73×54=3942", "58×64=3712"),
    @("63×98=6174", "13×71=923"),
    @("87×54=4698", "93×88=8184"),
    @("71×79=5609", "96×24=2304")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
